$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "passive income"
$ws.Range("B7").Value = "passive.income.nadi.myfirstdrawermenuproject"

$ws.Rows.Item(7).RowHeight = 23.85

$ws.Range("B7").Select()
